$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Prerequisites"
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"
$ws.Range("G1").Value = "Terms Typically Offered"

$ws.Range("C2").Value = "BCHEM/CHEM majors only."
$ws.Range("D2").Value = "NA"
$ws.Range("E2").Value = "NA"
$ws.Range("F2").Value = "NA"
$ws.Range("G2").Value = "F"

$ws.Range("C3").Value = "MATH 96; or MATH 115; or appropriate Math Placement Level."
$ws.Range("D3").Value = "NA"
$ws.Range("E3").Value = "NA"
$ws.Range("F3").Value = "NA"
$ws.Range("G3").Value = "F,W,SP,SU"

$ws.Range("C4").Value = "MATH 118."
$ws.Range("D4").Value = "NA"
$ws.Range("E4").Value = "NA"
$ws.Range("F4").Value = "High school chemistry or equivalent."
$ws.Range("G4").Value = "F,W,SP,SU "

$ws.Range("C5").Value = "CHEM 124, or AP Chemistry score of 5."
$ws.Range("D5").Value = "NA"
$ws.Range("E5").Value = "NA"
$ws.Range("F5").Value = "NA"
$ws.Range("G5").Value = "F, W, SP"

$ws.Range("C6").Value = "CHEM 125 with a grade of C- or better or consent of instructor."
$ws.Range("D6").Value = "NA"
$ws.Range("E6").Value = "NA"
$ws.Range("F6").Value = "NA"
$ws.Range("G6").Value = "SP"

$ws.Range("C7").Value = "MATH 118 or MATH 330."
$ws.Range("D7").Value = "NA"
$ws.Range("E7").Value = "NA"
$ws.Range("F7").Value = "High school chemistry or equivalent."
$ws.Range("G7").Value = "F, W "

$ws.Range("C8").Value = "CHEM 127 or AP Chemistry score of 5."
$ws.Range("D8").Value = "NA"
$ws.Range("E8").Value = "NA"
$ws.Range("F8").Value = "NA"
$ws.Range("G8").Value = "W, SP"

$ws.Range("C9").Value = "CHEM 128."
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "NA"
$ws.Range("F9").Value = "NA"
$ws.Range("G9").Value = "F, SP"

$ws.Range("C10").Value = "CHEM 111, CHEM 124, or CHEM 127 and consent of department chair."
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = "NA"
$ws.Range("F10").Value = "NA"
$ws.Range("G10").Value = "TBD"

$ws.Range("C11").Value = "Consent of instructor."
$ws.Range("D11").Value = "NA"
$ws.Range("E11").Value = "NA"
$ws.Range("F11").Value = "NA"
$ws.Range("G11").Value = "F,W,SP,SU"

$ws.Range("C12").Value = "Completion of a course with a BIO, BOT or MCRO prefix and a course with a CHEM prefix."
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("G12").Value = "F, W"

$ws.Range("C13").Value = "CHEM 126."
$ws.Range("D13").Value = "NA"
$ws.Range("E13").Value = "NA"
$ws.Range("F13").Value = "NA"
$ws.Range("G13").Value = "W, SP"

$ws.Range("C14").Value = "CHEM 111, CHEM 124 or CHEM 127."
$ws.Range("D14").Value = "NA"
$ws.Range("E14").Value = "NA"
$ws.Range("F14").Value = "NA"
$ws.Range("G14").Value = "TBD"

$ws.Range("C15").Value = "CHEM 126 or CHEM 129 with a grade of C- or better or consent of instructor."
$ws.Range("D15").Value = "NA"
$ws.Range("E15").Value = "NA"
$ws.Range("F15").Value = "NA"
$ws.Range("G15").Value = "F, W"

$ws.Range("C16").Value = "CHEM 216 with a grade of C- or better or consent of instructor."
$ws.Range("D16").Value = "CHEM 221 for Chemistry and Biochemistry majors; or CHEM 220 for non-Chemistry and non-Biochemistry majors."
$ws.Range("E16").Value = "NA"
$ws.Range("F16").Value = "NA"
$ws.Range("G16").Value = "W, SP "

$ws.Range("C17").Value = "CHEM 217 with a grade of C- or better or consent of instructor."
$ws.Range("D17").Value = "CHEM 324 for Chemistry and Biochemistry majors; or CHEM 223 for non-Chemistry and non-Biochemistry majors."
$ws.Range("E17").Value = "NA"
$ws.Range("F17").Value = "NA"
$ws.Range("G17").Value = "F, SP "

$ws.Range("C18").Value = "NA"
$ws.Range("D18").Value = "CHEM 217."
$ws.Range("E18").Value = "NA"
$ws.Range("F18").Value = "NA"
$ws.Range("G18").Value = "W, SP"

$ws.Range("C19").Value = "major in Chemistry or Biochemistry."
$ws.Range("D19").Value = "CHEM 217."
$ws.Range("E19").Value = "NA"
$ws.Range("F19").Value = "NA"
$ws.Range("G19").Value = "W, SP "

$ws.Range("C20").Value = "NA"
$ws.Range("D20").Value = "CHEM 218."
$ws.Range("E20").Value = "NA"
$ws.Range("F20").Value = "NA"
$ws.Range("G20").Value = "F, SP"

$ws.Range("C21").Value = "CHEM 126 or 129."
$ws.Range("D21").Value = "NA"
$ws.Range("E21").Value = "NA"
$ws.Range("F21").Value = "NA"
$ws.Range("G21").Value = "TBD"

$ws.Range("C22").Value = "CHEM 111, CHEM 124 or CHEM 127."
$ws.Range("D22").Value = "NA"
$ws.Range("E22").Value = "NA"
$ws.Range("F22").Value = "NA"
$ws.Range("G22").Value = "TBD"

$ws.Range("C23").Value = "Open to undergraduate students and consent of instructor."
$ws.Range("D23").Value = "NA"
$ws.Range("E23").Value = "NA"
$ws.Range("F23").Value = "NA"
$ws.Range("G23").Value = "TBD"

$ws.Range("C24").Value = "CHEM 216 or CHEM 312."
$ws.Range("D24").Value = "NA"
$ws.Range("E24").Value = "NA"
$ws.Range("F24").Value = "NA"
$ws.Range("G24").Value = "SP"

$ws.Range("C25").Value = "CHEM 203 and CHEM 218."
$ws.Range("D25").Value = "NA"
$ws.Range("E25").Value = "NA"
$ws.Range("F25").Value = "NA"
$ws.Range("G25").Value = "F, W, SP"

$ws.Range("C26").Value = "Junior standing; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; completion of GE Area B2; and one of the following CHEM 110, CHEM 124, or CHEM 127."
$ws.Range("D26").Value = "NA"
$ws.Range("E26").Value = "NA"
$ws.Range("F26").Value = "NA"
$ws.Range("G26").Value = "W"

$ws.Range("C27").Value = "CHEM 125 or CHEM 128."
$ws.Range("D27").Value = "NA"
$ws.Range("E27").Value = "NA"
$ws.Range("F27").Value = "NA"
$ws.Range("G27").Value = "F,W,SP,SU"

$ws.Range("C28").Value = "CHEM 212, CHEM 216, CHEM 312, or CHEM 316."
$ws.Range("D28").Value = "NA"
$ws.Range("E28").Value = "NA"
$ws.Range("F28").Value = "NA"
$ws.Range("G28").Value = "F, W, SP"

$ws.Range("C29").Value = "major in Chemistry or Biochemistry."
$ws.Range("D29").Value = "CHEM 218."
$ws.Range("E29").Value = "NA"
$ws.Range("F29").Value = "NA"
$ws.Range("G29").Value = "F, SP "

$ws.Range("C30").Value = "CHEM 126 or 129."
$ws.Range("D30").Value = "NA"
$ws.Range("E30").Value = "NA"
$ws.Range("F30").Value = "NA"
$ws.Range("G30").Value = "F, SP, SU"

$ws.Range("C31").Value = "CHEM 126 or 129; and CHEM 212 or CHEM 312; or CHEM 216 or CHEM 316."
$ws.Range("D31").Value = "NA"
$ws.Range("E31").Value = "NA"
$ws.Range("F31").Value = "NA"
$ws.Range("G31").Value = "TBD"

$ws.Range("C32").Value = "Junior standing; completion of GE Area A with grades of C- or better; completion of GE Area B1 with a grade of C- or better in at least one of the courses; BIO or MCRO course in GE Area B2; CHEM course in GE Area B3; and completion of GE Area B4."
$ws.Range("D32").Value = "NA"
$ws.Range("E32").Value = "NA"
$ws.Range("F32").Value = "NA"
$ws.Range("G32").Value = "F, SP"

$ws.Range("C33").Value = "CHEM 126 or CHEM 129; MATH 143; PHYS 122 or PHYS 132."
$ws.Range("D33").Value = "NA"
$ws.Range("E33").Value = "NA"
$ws.Range("F33").Value = "NA"
$ws.Range("G33").Value = "F, W"

$ws.Range("C34").Value = "CHEM 351."
$ws.Range("D34").Value = "NA"
$ws.Range("E34").Value = "NA"
$ws.Range("F34").Value = "NA"
$ws.Range("G34").Value = "W, SP"

$ws.Range("C35").Value = "CHEM 352."
$ws.Range("D35").Value = "NA"
$ws.Range("E35").Value = "NA"
$ws.Range("F35").Value = "NA"
$ws.Range("G35").Value = "F, SP"

$ws.Range("C36").Value = "CHEM 231/331."
$ws.Range("D36").Value = "CHEM 352."
$ws.Range("E36").Value = "NA"
$ws.Range("F36").Value = "NA"
$ws.Range("G36").Value = "F, W, SP "

$ws.Range("C37").Value = "NA"
$ws.Range("D37").Value = "CHEM 353."
$ws.Range("E37").Value = "NA"
$ws.Range("F37").Value = "NA"
$ws.Range("G37").Value = "F, SP"

$ws.Range("C38").Value = "CHEM 217 or CHEM 317; and BIO 161."
$ws.Range("D38").Value = "NA"
$ws.Range("E38").Value = "NA"
$ws.Range("F38").Value = "CHEM 231/331."
$ws.Range("G38").Value = "F, W, SP "

$ws.Range("C39").Value = "CHEM 371."
$ws.Range("D39").Value = "NA"
$ws.Range("E39").Value = "NA"
$ws.Range("F39").Value = "NA"
$ws.Range("G39").Value = "F, SP"

$ws.Range("C40").Value = "CHEM 371."
$ws.Range("D40").Value = "NA"
$ws.Range("E40").Value = "NA"
$ws.Range("F40").Value = "NA"
$ws.Range("G40").Value = "W, SP"

$ws.Range("C41").Value = "CHEM 313 or CHEM 371."
$ws.Range("D41").Value = "NA"
$ws.Range("E41").Value = "NA"
$ws.Range("F41").Value = "NA"
$ws.Range("G41").Value = "W, SP"

$ws.Range("C42").Value = "Junior standing and consent of department chair."
$ws.Range("D42").Value = "NA"
$ws.Range("E42").Value = "NA"
$ws.Range("F42").Value = "NA"
$ws.Range("G42").Value = "TBD"

$ws.Range("C43").Value = "Consent of instructor."
$ws.Range("D43").Value = "NA"
$ws.Range("E43").Value = "NA"
$ws.Range("F43").Value = "NA"
$ws.Range("G43").Value = "F,W,SP,SU"

$ws.Range("C44").Value = "CHEM 303 and CHEM 352."
$ws.Range("D44").Value = "NA"
$ws.Range("E44").Value = "NA"
$ws.Range("F44").Value = "NA"
$ws.Range("G44").Value = "F, W, SP"

$ws.Range("C45").Value = "NA"
$ws.Range("D45").Value = "CHEM 353."
$ws.Range("E45").Value = "NA"
$ws.Range("F45").Value = "NA"
$ws.Range("G45").Value = "TBD"

$ws.Range("C46").Value = "CHEM 218 or CHEM 318."
$ws.Range("D46").Value = "NA"
$ws.Range("E46").Value = "NA"
$ws.Range("F46").Value = "NA"
$ws.Range("G46").Value = "W"

$ws.Range("C47").Value = "BIO 161 and CHEM 217."
$ws.Range("D47").Value = "NA"
$ws.Range("E47").Value = "NA"
$ws.Range("F47").Value = "NA"
$ws.Range("G47").Value = "F"

$ws.Range("C48").Value = "CHEM 218 or CHEM 318; CHEM 313 or CHEM 371."
$ws.Range("D48").Value = "NA"
$ws.Range("E48").Value = "NA"
$ws.Range("F48").Value = "NA"
$ws.Range("G48").Value = "TBD"

$ws.Range("C49").Value = "CHEM 218/318."
$ws.Range("D49").Value = "NA"
$ws.Range("E49").Value = "NA"
$ws.Range("F49").Value = "NA"
$ws.Range("G49").Value = "SP"

$ws.Range("C50").Value = "CHEM 371."
$ws.Range("D50").Value = "NA"
$ws.Range("E50").Value = "NA"
$ws.Range("F50").Value = "NA"
$ws.Range("G50").Value = "TBD"

$ws.Range("C51").Value = "CHEM 231/331, CHEM 354."
$ws.Range("D51").Value = "NA"
$ws.Range("E51").Value = "NA"
$ws.Range("F51").Value = "CHEM 353."
$ws.Range("G51").Value = "W "

$ws.Range("C52").Value = "Junior standing; BIO 161 or BIO 303."
$ws.Range("D52").Value = "NA"
$ws.Range("E52").Value = "NA"
$ws.Range("F52").Value = "BIO 302 or BIO 303 or BIO 351 or CHEM 373."
$ws.Range("G52").Value = "F, SP "

$ws.Range("C53").Value = "CHEM 212/312 or CHEM 216/316."
$ws.Range("D53").Value = "NA"
$ws.Range("E53").Value = "NA"
$ws.Range("F53").Value = "NA"
$ws.Range("G53").Value = "F"

$ws.Range("C54").Value = "CHEM 217/317 and CHEM 444."
$ws.Range("D54").Value = "NA"
$ws.Range("E54").Value = "NA"
$ws.Range("F54").Value = "NA"
$ws.Range("G54").Value = "W"

$ws.Range("C55").Value = "CHEM 125 or CHEM 128; CHEM 351, MATE 380, or ME 302."
$ws.Range("D55").Value = "NA"
$ws.Range("E55").Value = "NA"
$ws.Range("F55").Value = "NA"
$ws.Range("G55").Value = "SP"

$ws.Range("C56").Value = "NA"
$ws.Range("D56").Value = "CHEM 444."
$ws.Range("E56").Value = "NA"
$ws.Range("F56").Value = "NA"
$ws.Range("G56").Value = "F"

$ws.Range("C57").Value = "CHEM 447."
$ws.Range("D57").Value = "CHEM 445."
$ws.Range("E57").Value = "NA"
$ws.Range("F57").Value = "NA"
$ws.Range("G57").Value = "W "

$ws.Range("C58").Value = "CHEM 444."
$ws.Range("D58").Value = "NA"
$ws.Range("E58").Value = "NA"
$ws.Range("F58").Value = "NA"
$ws.Range("G58").Value = "F,W,SP,SU"

$ws.Range("C59").Value = "CHEM 444 or CHEM 544."
$ws.Range("D59").Value = "NA"
$ws.Range("E59").Value = "NA"
$ws.Range("F59").Value = "NA"
$ws.Range("G59").Value = "SP"

$ws.Range("C60").Value = "CHEM 447 or CHEM 547."
$ws.Range("D60").Value = "CHEM 450."
$ws.Range("E60").Value = "NA"
$ws.Range("F60").Value = "CHEM 445 or CHEM 545; CHEM 448 or CHEM 548; CHEM 446."
$ws.Range("G60").Value = "SP  "

$ws.Range("C61").Value = "CHEM 212 or CHEM 216 or CHEM 312 or CHEM 316; CHEM 351 or MATE 380; or graduate standing."
$ws.Range("D61").Value = "NA"
$ws.Range("E61").Value = "NA"
$ws.Range("F61").Value = "NA"
$ws.Range("G61").Value = "SP"

$ws.Range("C62").Value = "CHEM 324."
$ws.Range("D62").Value = "NA"
$ws.Range("E62").Value = "NA"
$ws.Range("F62").Value = "NA"
$ws.Range("G62").Value = "SP"

$ws.Range("C63").Value = "NA"
$ws.Range("D63").Value = "CHEM 218 or CHEM 318 and junior standing."
$ws.Range("E63").Value = "NA"
$ws.Range("F63").Value = "NA"
$ws.Range("G63").Value = "F,W,SP,SU"

$ws.Range("C64").Value = "Consent of instructor."
$ws.Range("D64").Value = "NA"
$ws.Range("E64").Value = "NA"
$ws.Range("F64").Value = "NA"
$ws.Range("G64").Value = "F,W,SP,SU"

$ws.Range("C65").Value = "Junior standing and consent of instructor."
$ws.Range("D65").Value = "NA"
$ws.Range("E65").Value = "NA"
$ws.Range("F65").Value = "NA"
$ws.Range("G65").Value = "F,W,SP,SU"

$ws.Range("C66").Value = "Junior standing, CHEM 231/331 (or permission of instructor), evidence of satisfactory preparation in chemistry; department chair approval required."
$ws.Range("D66").Value = "NA"
$ws.Range("E66").Value = "NA"
$ws.Range("F66").Value = "NA"
$ws.Range("G66").Value = "F, W, SP"

$ws.Range("C67").Value = "Junior standing and consent of instructor."
$ws.Range("D67").Value = "NA"
$ws.Range("E67").Value = "NA"
$ws.Range("F67").Value = "NA"
$ws.Range("G67").Value = "F, W, SP"

$ws.Range("C68").Value = "CHEM 351, CHEM 217 or CHEM 317."
$ws.Range("D68").Value = "NA"
$ws.Range("E68").Value = "NA"
$ws.Range("F68").Value = "NA"
$ws.Range("G68").Value = "F"

$ws.Range("C69").Value = "Consent of instructor."
$ws.Range("D69").Value = "NA"
$ws.Range("E69").Value = "NA"
$ws.Range("F69").Value = "NA"
$ws.Range("G69").Value = "TBD"

$ws.Range("C70").Value = "CHEM 371."
$ws.Range("D70").Value = "NA"
$ws.Range("E70").Value = "NA"
$ws.Range("F70").Value = "NA"
$ws.Range("G70").Value = "F, W, SP"

$ws.Range("C71").Value = "BIO 161, and grade of C- or better in BIO 351 or CHEM 373 or consent of instructor."
$ws.Range("D71").Value = "NA"
$ws.Range("E71").Value = "NA"
$ws.Range("F71").Value = "NA"
$ws.Range("G71").Value = "F, W, SP"

$ws.Range("C72").Value = "BIO/CHEM 475; CHEM 313 or CHEM 371, or graduate standing in Biological Sciences."
$ws.Range("D72").Value = "NA"
$ws.Range("E72").Value = "NA"
$ws.Range("F72").Value = "NA"
$ws.Range("G72").Value = "SP"

$ws.Range("C73").Value = "CHEM 218 or CHEM 318."
$ws.Range("D73").Value = "NA"
$ws.Range("E73").Value = "NA"
$ws.Range("F73").Value = "NA"
$ws.Range("G73").Value = "TBD"

$ws.Range("C74").Value = "CHEM 352 and CHEM 231/331."
$ws.Range("D74").Value = "NA"
$ws.Range("E74").Value = "NA"
$ws.Range("F74").Value = "NA"
$ws.Range("G74").Value = "F"

$ws.Range("C75").Value = "NA"
$ws.Range("D75").Value = "CHEM 481."
$ws.Range("E75").Value = "NA"
$ws.Range("F75").Value = "NA"
$ws.Range("G75").Value = "F"

$ws.Range("C76").Value = "Sophomore standing and consent of instructor."
$ws.Range("D76").Value = "NA"
$ws.Range("E76").Value = "NA"
$ws.Range("F76").Value = "NA"
$ws.Range("G76").Value = "F,W,SP,SU"

$ws.Range("C77").Value = "Sophomore standing and consent of instructor."
$ws.Range("D77").Value = "NA"
$ws.Range("E77").Value = "NA"
$ws.Range("F77").Value = "NA"
$ws.Range("G77").Value = "F,W,SP,SU"

$ws.Range("C78").Value = "Graduate standing and consent of department chair."
$ws.Range("D78").Value = "NA"
$ws.Range("E78").Value = "NA"
$ws.Range("F78").Value = "NA"
$ws.Range("G78").Value = "F,W,SP,SU"

$ws.Range("C79").Value = "CHEM 212/312 or CHEM 216/316 or equivalent; CHEM 351 or equivalent."
$ws.Range("D79").Value = "NA"
$ws.Range("E79").Value = "NA"
$ws.Range("F79").Value = "NA"
$ws.Range("G79").Value = "F"

$ws.Range("C80").Value = "CHEM 544."
$ws.Range("D80").Value = "NA"
$ws.Range("E80").Value = "NA"
$ws.Range("F80").Value = "NA"
$ws.Range("G80").Value = "W"

$ws.Range("C81").Value = "NA"
$ws.Range("D81").Value = "CHEM 544."
$ws.Range("E81").Value = "NA"
$ws.Range("F81").Value = "NA"
$ws.Range("G81").Value = "F"

$ws.Range("C82").Value = "CHEM 547."
$ws.Range("D82").Value = "CHEM 545."
$ws.Range("E82").Value = "NA"
$ws.Range("F82").Value = "NA"
$ws.Range("G82").Value = "W "

$ws.Range("C83").Value = "CHEM 444 or CHEM 544."
$ws.Range("D83").Value = "NA"
$ws.Range("E83").Value = "NA"
$ws.Range("F83").Value = "NA"
$ws.Range("G83").Value = "SP"

$ws.Range("C84").Value = "NA"
$ws.Range("D84").Value = "CHEM 550."
$ws.Range("E84").Value = "NA"
$ws.Range("F84").Value = "NA"
$ws.Range("G84").Value = "SP"

$ws.Range("C85").Value = "Graduate standing or consent of instructor."
$ws.Range("D85").Value = "NA"
$ws.Range("E85").Value = "NA"
$ws.Range("F85").Value = "NA"
$ws.Range("G85").Value = "TBD"

$ws.Range("C86").Value = "Graduate standing or consent of instructor."
$ws.Range("D86").Value = "NA"
$ws.Range("E86").Value = "NA"
$ws.Range("F86").Value = "NA"
$ws.Range("G86").Value = "TBD"

$ws.Range("C87").Value = "Graduate standing in the Polymers and Coatings program or consent of instructor."
$ws.Range("D87").Value = "NA"
$ws.Range("E87").Value = "NA"
$ws.Range("F87").Value = "NA"
$ws.Range("G87").Value = "F, W, SP"

$ws.Range("C88").Value = "CHEM 545, CHEM 547, CHEM 548, CHEM 550, CHEM 551."
$ws.Range("D88").Value = "NA"
$ws.Range("E88").Value = "NA"
$ws.Range("F88").Value = "NA"
$ws.Range("G88").Value = "F,W,SP,SU"

$ws.Range("C89").Value = "CHEM 545, CHEM 547, CHEM 548, CHEM 550, CHEM 551."
$ws.Range("D89").Value = "NA"
$ws.Range("E89").Value = "NA"
$ws.Range("F89").Value = "NA"
$ws.Range("G89").Value = "F,W,SP,SU"
